$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    if ($text -eq "") {
        $cell.Value = "'"
    } else {
        $cell.Value = $text
    }
    $cell.Style = "Normal"
}

# Row 6
Set-TextCell 6 1 'Kiwi Kitchens Raw Freeze Dried Cat Treats Salmon Recipe 30g | Catch.com.au'
Set-TextCell 6 2 '$6.90'
Set-TextCell 6 3 'https://www.catch.com.au/product/kiwi-kitchens-raw-freeze-dried-cat-treats-salmon-recipe-30g-25308308/?sid=Kiwi%20Kitchens%20Raw%20Freeze%20Dried%20Cat%20Treats%20Tuna%20Recipe%2030g&st=32&srtrev=sj-j5nztvbthwy0bjea4vf5q8.click&pid=25308308&sp=2&oid=90790459'
$ws.Cells.Item(6, 4).Value = 6.9

# Row 7
Set-TextCell 7 1 'OsteVit-D One-A-Day Vitamin D3 250 Tabs | Catch.com.au'
Set-TextCell 7 2 '$12.99'
Set-TextCell 7 3 'https://www.catch.com.au/product/ostevit-d-one-a-day-vitamin-d3-250-tabs-26174042/?sid=Grocery%20%26%20Liquor%20%3E%20Health%20%26%20Wellbeing%20%3E%20Vitamins&st=15&srtrev=sj-32jmfaklbejowe1epxavf2.click&pid=26174042&sp=26&oid=97836753'
$ws.Cells.Item(7, 4).Value = 12.99

# Row 8
Set-TextCell 8 1 'BeePower Manuka Honey 500+ MGO Lozenges Lemon 40pk | Catch.com.au'
Set-TextCell 8 2 '$16.19'
Set-TextCell 8 3 'https://www.catch.com.au/product/beepower-manuka-honey-500-mgo-lozenges-lemon-40pk-6119470/?sid=Bee%20Power%20Manuka%20Honey%20Lemon%20Lozenges%20190%20g&st=32&srtrev=sj-dz00y66djmu4900pcp04ht.click&pid=6119470&sp=1&oid=30391320'
$ws.Cells.Item(8, 4).Value = 14.39

# Row 9
Set-TextCell 9 1 'Nature''s Way Restore Probiotic Daily Health 28 Caps | Catch.com.au'
Set-TextCell 9 2 '$16.99'
Set-TextCell 9 3 'https://www.catch.com.au/product/natures-way-restore-probiotic-daily-health-28-caps-7728490/?sid=Grocery%20%26%20Liquor%20%3E%20Health%20%26%20Wellbeing%20%3E%20Vitamins&st=15&srtrev=sj-6ak5i7pvzaylbfqzjqwaud.click&pid=7728490&sp=15&oid=39164780'
$ws.Cells.Item(9, 4).Value = 16.99

# Row 10
Set-TextCell 10 1 'Morlife Organic Barley Grass Powder 700g / 350 Serves | Catch.com.au'
Set-TextCell 10 2 '$44.95'
Set-TextCell 10 3 'https://www.catch.com.au/product/morlife-organic-barley-grass-powder-700g-350-serves-417089/?sid=Grocery%20%26%20Liquor%20%3E%20Health%20%26%20Wellbeing%20%3E%20Diet%20%26%20Slimming&st=15&srtrev=sj-czn49tg62b20v2utysy9ki.click&pid=417089&sp=4&oid=462533'
$ws.Cells.Item(10, 4).Value = 44.95

# Row 11
Set-TextCell 11 1 '10 x Amazonia RAW Protein Bars Triple Choc Brownie 40g | Catch.com.au'
Set-TextCell 11 2 '$35'
Set-TextCell 11 3 'https://www.catch.com.au/product/10-x-amazonia-raw-protein-bars-triple-choc-brownie-40g-7730920/?sid=Grocery%20%26%20Liquor%20%3E%20Health%20%26%20Wellbeing%20%3E%20Diet%20%26%20Slimming&st=15&srtrev=sj-d789l8mhvrmo6j3ogjafly.click&pid=7730920&sp=14&oid=39174231'
$ws.Cells.Item(11, 4).Value = 35

# Row 12
Set-TextCell 12 1 'Piksters Interdental Brushes 40pk - Size 1 | Catch.com.au'
Set-TextCell 12 2 '$15.53'
Set-TextCell 12 3 'https://www.catch.com.au/product/piksters-interdental-brushes-40pk-size-1-6451946/?sid=Beauty%20%3E%20Dental%20%26%20Oral%20Care%20%3E%20Dental%20Floss%20%26%20Picks&st=15&srtrev=sj-c20sj0wcx4c5wub79fwjeg.click&pid=6451946&sp=7&oid=32153160'
$ws.Cells.Item(12, 4).Value = 15.53

# Row 13
Set-TextCell 13 1 'Joseph Joseph 4-Piece Nested Glass Storage Container Set | Catch.com.au'
Set-TextCell 13 2 ''
Set-TextCell 13 3 'https://www.catch.com.au/product/joseph-joseph-4-piece-nested-glass-storage-container-set-1911952/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen%20%3E%20Food%20Storage&st=15&srtrev=sj-pqadxjaq4zlnovpbblxvyi.click&pid=1911952&sp=16&oid=11989745'
$ws.Cells.Item(13, 4).Value = 39.2

# Row 14
Set-TextCell 14 1 'Avanti 38x27.5cm Potato Bag | Catch.com.au'
Set-TextCell 14 2 '$18.95'
Set-TextCell 14 3 'https://www.catch.com.au/product/avanti-38x27-5cm-potato-bag-2311534/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen%20%3E%20Food%20Storage&st=15&srtrev=sj-6yezpvrjxrceocdiwerew8.click&pid=2311534&sp=26&oid=75509167'
$ws.Cells.Item(14, 4).Value = 6.69

# Row 15
Set-TextCell 15 1 'Bakemaster 25cm Square Proving Basket | Catch.com.au'
Set-TextCell 15 2 '$12.95'
Set-TextCell 15 3 'https://www.catch.com.au/product/bakemaster-25cm-square-proving-basket-15086436/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen%20%3E%20Food%20Storage&st=15&srtrev=sj-rud6qpy8iwtc87rrfp4aus.click&pid=15086436&sp=95&oid=75509271'
$ws.Cells.Item(15, 4).Value = 12.95

# Row 16
Set-TextCell 16 1 'Maxwell & Williams 6-Cup BakerMaker Non-Stick Large Muffin Pan | Catch.com.au'
Set-TextCell 16 2 '$10.73'
Set-TextCell 16 3 'https://www.catch.com.au/product/maxwell-williams-6-cup-bakermaker-non-stick-large-muffin-pan-24879640/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen%20%3E%20Baking&st=15&srtrev=sj-btsljf1lfkv2ocf0yq3d4b.click&pid=24879640&sp=26&oid=87338530'
$ws.Cells.Item(16, 4).Value = 10.73

# Row 17
Set-TextCell 17 1 'Salt & Pepper 12-Piece Napier Dinner Set - Blue | Catch.com.au'
Set-TextCell 17 2 '$159.95'
Set-TextCell 17 3 'https://www.catch.com.au/product/salt-pepper-12-piece-napier-dinner-set-blue-11801874/?sid=Home%20%26%20Kitchen%20%3E%20Dining%20%3E%20Dinnerware%20%26%20Serveware&st=15&srtrev=sj-8dh3vrcme4ydu8egzc1rr0.click&pid=11801874&sp=1&oid=50160723'
$ws.Cells.Item(17, 4).Value = 35

# Row 18
Set-TextCell 18 1 'KitchenAid Soft Touch Can Opener | Catch.com.au'
Set-TextCell 18 2 '$14.95'
Set-TextCell 18 3 'https://www.catch.com.au/product/kitchenaid-soft-touch-can-opener-5027792/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen&st=15&srtrev=sj-h1anpw2bj8z4ficsigz3kl.click&pid=5027792&sp=23&oid=30079340'
$ws.Cells.Item(18, 4).Value = 14.95

# Row 19
Set-TextCell 19 1 'Maxwell & Williams 5-Piece Cocktail & Co. Lexington Hammered Bar Tool Set | Catch.com.au'
Set-TextCell 19 2 '$39'
Set-TextCell 19 3 'https://www.catch.com.au/product/maxwell-williams-5-piece-cocktail-co-lexington-hammered-bar-tool-set-24879477/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen&st=15&srtrev=sj-lhzqnuy6g22auci1hachod.click&pid=24879477&sp=34&oid=87337786'
$ws.Cells.Item(19, 4).Value = 39

# Row 20
Set-TextCell 20 1 'BarCraft Twist Action Wine Bottle Foil Cutter | Catch.com.au'
Set-TextCell 20 2 '$4'
Set-TextCell 20 3 'https://www.catch.com.au/product/barcraft-twist-action-wine-bottle-foil-cutter-6045358/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen&st=15&srtrev=sj-7dpx2ioq6qwdfqbvlusqjn.click&pid=6045358&sp=44&oid=31733262'
$ws.Cells.Item(20, 4).Value = 4

# Row 21
Set-TextCell 21 1 'Joseph Joseph Push & Tear Kitchen Roll/Paper Towel Stand | Catch.com.au'
Set-TextCell 21 2 '$49.95'
Set-TextCell 21 3 'https://www.catch.com.au/product/joseph-joseph-push-tear-kitchen-roll-paper-towel-stand-3904093/?sid=Home%20%26%20Kitchen%20%3E%20Kitchen&st=15&srtrev=sj-pnw3v1vp6p6w2t7xcitwjv.click&pid=3904093&sp=61&oid=29945864'
$ws.Cells.Item(21, 4).Value = 34.95
